$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.242871403694153
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.275062084197998
$ws.Range("D1").Value = 1.292945742607117
$ws.Range("E1").Value = 0.9337817430496216
